$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
# A8 shared rich-text run: "Volume 31   Number  33" -> "...  34"
$ws.Range("A8").Characters(21, 2).Text = "34"
# C9 shared rich-text runs: week "... 8/12/2024  Through  8/18/2024" -> "...8/19/2024 ... 8/25/2024"
$ws.Range("C9").Characters(27, 9).Text = "8/19/2024"
$ws.Range("C9").Characters(47, 9).Text = "8/25/2024"

# --- Cells changing between numeric and shared-text ("0" / "***.*") representation ---
# Each target state (value + style) already exists verbatim elsewhere on the sheet, so the
# whole cell is copied from a matching donor cell to reproduce it exactly.
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("G15").Copy($ws.Range("C20"))
$ws.Range("F14").Copy($ws.Range("C25"))
$ws.Range("C17").Copy($ws.Range("C28"))
$ws.Range("D14").Copy($ws.Range("D29"))
$ws.Range("H22").Copy($ws.Range("E29"))
$ws.Range("F14").Copy($ws.Range("D30"))
$ws.Range("N22").Copy($ws.Range("E30"))

# --- Updated crime-statistics figures ---
$ws.Range("G14").Value = 1
$ws.Range("N14").Value = -81.081081081081
$ws.Range("L15").Value = -47.368421052631
$ws.Range("N15").Value = -83.606557377049
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -40
$ws.Range("F16").Value = 21
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = 40
$ws.Range("I16").Value = 154
$ws.Range("J16").Value = 133
$ws.Range("K16").Value = 15.78947368421
$ws.Range("L16").Value = 12.408759124087
$ws.Range("M16").Value = -19.791666666666
$ws.Range("N16").Value = -73.448275862069
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -37.5
$ws.Range("G17").Value = 41
$ws.Range("H17").Value = -34.146341463414
$ws.Range("I17").Value = 311
$ws.Range("J17").Value = 307
$ws.Range("K17").Value = 1.302931596091
$ws.Range("L17").Value = -10.632183908046
$ws.Range("M17").Value = 42.009132420091
$ws.Range("N17").Value = -51.857585139318
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -40
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -28.571428571428
$ws.Range("I18").Value = 83
$ws.Range("J18").Value = 122
$ws.Range("K18").Value = -31.967213114754
$ws.Range("L18").Value = -37.593984962406
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = -86.846275752773
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -55.555555555555
$ws.Range("G19").Value = 45
$ws.Range("H19").Value = -48.888888888888
$ws.Range("I19").Value = 222
$ws.Range("J19").Value = 283
$ws.Range("K19").Value = -21.554770318021
$ws.Range("L19").Value = -2.202643171806
$ws.Range("M19").Value = 14.432989690721
$ws.Range("N19").Value = -9.38775510204
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = -87.5
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 25
$ws.Range("H20").Value = -88
$ws.Range("I20").Value = 58
$ws.Range("J20").Value = 83
$ws.Range("K20").Value = -30.12048192771
$ws.Range("L20").Value = 11.538461538461
$ws.Range("M20").Value = 70.588235294117
$ws.Range("N20").Value = -65.882352941176
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 35
$ws.Range("E21").Value = -54.285714285714
$ws.Range("F21").Value = 84
$ws.Range("G21").Value = 142
$ws.Range("H21").Value = -40.845070422535
$ws.Range("I21").Value = 845
$ws.Range("J21").Value = 954
$ws.Range("K21").Value = -11.425576519916
$ws.Range("L21").Value = -8.351409978308
$ws.Range("M21").Value = 13.270777479892
$ws.Range("N21").Value = -64.345991561181
$ws.Range("L22").Value = -53.333333333333
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -33.333333333333
$ws.Range("F23").Value = 15
$ws.Range("G23").Value = 21
$ws.Range("H23").Value = -28.571428571428
$ws.Range("I23").Value = 178
$ws.Range("J23").Value = 155
$ws.Range("K23").Value = 14.838709677419
$ws.Range("L23").Value = 23.611111111111
$ws.Range("M23").Value = 58.928571428571
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = -21.428571428571
$ws.Range("F24").Value = 55
$ws.Range("G24").Value = 78
$ws.Range("H24").Value = -29.487179487179
$ws.Range("I24").Value = 540
$ws.Range("J24").Value = 569
$ws.Range("K24").Value = -5.096660808435
$ws.Range("L24").Value = -4.761904761904
$ws.Range("M24").Value = 30.750605326876
$ws.Range("E25").Value = -100
$ws.Range("F25").Value = 9
$ws.Range("H25").Value = 12.5
$ws.Range("J25").Value = 113
$ws.Range("K25").Value = -10.619469026548
$ws.Range("L25").Value = -31.756756756756
$ws.Range("C26").Value = 17
$ws.Range("D26").Value = 9
$ws.Range("E26").Value = 88.888888888888
$ws.Range("F26").Value = 71
$ws.Range("G26").Value = 37
$ws.Range("H26").Value = 91.891891891891
$ws.Range("I26").Value = 476
$ws.Range("J26").Value = 386
$ws.Range("K26").Value = 23.316062176165
$ws.Range("L26").Value = 23.316062176165
$ws.Range("M26").Value = -17.073170731707
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -66.666666666666
$ws.Range("I27").Value = 16
$ws.Range("J27").Value = 33
$ws.Range("K27").Value = -51.515151515151
$ws.Range("L27").Value = -50
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 35
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = -36.363636363636
$ws.Range("G29").Value = 5
$ws.Range("L29").Value = -44.444444444444
$ws.Range("M29").Value = -48.275862068965
$ws.Range("N29").Value = -81.481481481481
$ws.Range("G30").Value = 4
$ws.Range("L30").Value = -47.826086956521
$ws.Range("M30").Value = -50
$ws.Range("N30").Value = -84
